$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = -1.420907706991966
$ws.Range("B2").Value = -4.112924357213253

$ws.Range("A3").Value = -0.5063990210886026
$ws.Range("B3").Value = 0.9844700692951912

$ws.Range("A4").Value = 1.01086908773511
$ws.Range("B4").Value = -3.127852703437672

$ws.Range("A5").Value = 0.6814715563138282
$ws.Range("B5").Value = 0.4320455422917935

$ws.Range("A6").Value = -0.7885271013357633
$ws.Range("B6").Value = -1.859278549900989

$ws.Range("A7").Value = -0.05860466833245193
$ws.Range("B7").Value = -0.9113797029307716

$ws.Range("A8").Value = 0.790215535803325
$ws.Range("B8").Value = 0.706972271541737

$ws.Range("A9").Value = 0.313359721177403
$ws.Range("B9").Value = 1.081627833807036

$ws.Range("A10").Value = -0.1723433803118873
$ws.Range("B10").Value = -2.154301349311564

$ws.Range("A11").Value = 0.3418316664756613
$ws.Range("B11").Value = -0.8406096601245452

$ws.Range("A12").Value = 0.5328414683185517
$ws.Range("B12").Value = 0.1656461003204004
